$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Norway sheet: duplicate "Hungary" (the current last sheet), rename,
# update the market-name / product-code cells, and insert the extra
# "MZXSDR240" row that Norway's panel list needs (between MZX64DR and
# MZXDR240).
# ---------------------------------------------------------------------
$hungary = $wb.Worksheets.Item("Hungary")

$hungary.Copy([System.Reflection.Missing]::Value, $hungary)
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"

# New row for MZXSDR240, inserted above the existing MZXDR240 row (row 15)
$norway.Rows.Item(15).Insert()
$norway.Range("A14").Copy()
$norway.Range("A15").PasteSpecial(-4122)
$norway.Range("A15").Value = "MZXSDR240"

# Product code first, then market name - matches the order the strings
# were appended to the shared-string table in the source edit.
$norway.Range("B4").Value = "NGC-2931/T3063"
$norway.Range("B2").Value = "NorwayMarket"

# Column widths specific to the new country sheets.
$norway.Columns.Item(1).ColumnWidth = 21.5546875
$norway.Columns.Item(2).ColumnWidth = 24.109375
$norway.Columns.Item(3).ColumnWidth = 12.33203125
$norway.Columns.Item(4).ColumnWidth = 12.77734375

# Re-fit rows 3 & 4 now that the columns are wide enough that the text no
# longer wraps onto a second line (row 5 keeps its taller height).
$norway.Rows.Item(3).AutoFit()
$norway.Rows.Item(4).AutoFit()

$norway.Activate()
$norway.Range("A15").Select()

# ---------------------------------------------------------------------
# Poland sheet: duplicate "Hungary" again (standard panel list, no extra
# row), rename, update market-name / product-code cells.
# ---------------------------------------------------------------------
$hungary.Copy([System.Reflection.Missing]::Value, $norway)
$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"

$poland.Range("B4").Value = "NGC-2920/T3106"
$poland.Range("B2").Value = "Poland Market"

$poland.Columns.Item(1).ColumnWidth = 21.5546875
$poland.Columns.Item(2).ColumnWidth = 24.109375
$poland.Columns.Item(3).ColumnWidth = 12.33203125
$poland.Columns.Item(4).ColumnWidth = 12.77734375

$poland.Rows.Item(3).AutoFit()
$poland.Rows.Item(4).AutoFit()

$poland.Activate()
$poland.Range("A15").Select()
